$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same: Main / Xpath / Value

# Values (column C) first, then Xpaths (column B) - matches the order the
# strings were originally authored in, so the shared-string table comes out
# in the same sequence as the target workbook.
$ws.Range("C2").Value = "[A-Z a-z].*"
$ws.Range("C3").Value = "[A-Z]{2}"
$ws.Range("C4").Value = "[0-9]{5}"
$ws.Range("C5").Value = "[0-9]{3}"
$ws.Range("C6").Value = "[A-Z]{1}"

$ws.Range("B2").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table[1]/CITY"
$ws.Range("B3").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table[1]/STATE"
$ws.Range("B4").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table[1]/ZIP"
$ws.Range("B5").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table[1]/AREA_CODE"
$ws.Range("B6").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table[1]/TIME_ZONE"

# Column widths (engine snaps ColumnWidth to a pixel grid on save, so the
# inputs below are chosen so the saved width lands as close as possible to
# the target stored widths of 88.5703125 / 10 characters)
$ws.Columns.Item(2).ColumnWidth = 87.65
$ws.Columns.Item(3).ColumnWidth = 9.17

# Selection
$ws.Range("E8").Select()

$wb.Save()
